$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: move 4 remaining from Week1 to Week2 (E5 now 2, F5 now 4)
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 4

# Row 12: move the 6 remaining from Week1 to Week2 (E12 now blank, F12 now 6)
$ws.Range("E12").ClearContents()
$ws.Range("F12").Value = 6

# Update the active selection to F12 (was F13)
$ws.Range("F12").Select()
